$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-apply the existing date number-format style (same style already used
# by A2/A3/A7/A9/A17) to the new cells BEFORE the formulas are entered, so
# the engine doesn't auto-create a brand-new "General -> date" style and
# instead simply reuses the current one.
$ws.Range("A17").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add the new formulas in rows 18 and 19, column A
$ws.Range("A18").Formula = "=DATE(2000, 1, 1)"
$ws.Range("A19").Formula = '=DATE("2000",1, 1)'

# Set the selection to match: A4:D4, active cell A4
$ws.Range("A4:D4").Select()
